$wb = $excel.ActiveWorkbook

# --- Step 1: Fix "ParkEiffel" spelling to "ParkEifel" and rearrange parks alphabetically. ---
# This relabels rows 5,6,7,10 (the Park categories stored at those rows) on every sheet;
# rows 8 (ParkJasmund) and 9 (ParkKellerwald) already sort correctly and keep their text/values.
$parkLabelByRow = @{ 5 = "ParkEifel"; 6 = "ParkHainich"; 7 = "ParkHunsrueck"; 10 = "ParkSaechs_Schw" }

for ($s = 1; $s -le 8; $s++) {
    $ws = $wb.Worksheets.Item($s)
    foreach ($row in $parkLabelByRow.Keys) {
        $ws.Cells.Item($row, 1).Value = $parkLabelByRow[$row]
    }
}

# --- Sheet 1 (API): update coefficient table values ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 0.2467
$ws.Range("C2").Value = 85816.1197
$ws.Range("B5").Value = 1.7447
$ws.Range("C5").Value = 85816.1197
$ws.Range("B6").Value = 1.8242
$ws.Range("C6").Value = 85816.1197
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = 0
$ws.Range("D7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = 0
$ws.Range("D8").ClearContents()
$ws.Range("B9").Value = 1.3494
$ws.Range("C9").Value = 85816.1197
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 0
$ws.Range("D10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 0
$ws.Range("D11").ClearContents()

# --- Sheet 2 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B5").Value = -0.2778
$ws.Range("C5").Value = 0.1692
$ws.Range("D5").Value = 0.1006
$ws.Range("B6").Value = 0.2237
$ws.Range("C6").Value = 0.1549
$ws.Range("D6").Value = 0.1486
$ws.Range("B7").Value = -0.5301
$ws.Range("C7").Value = 0.2023
$ws.Range("D7").Value = 0.0088
$ws.Range("B10").Value = 0.241
$ws.Range("C10").Value = 0.147
$ws.Range("D10").Value = 0.1011

# --- Sheet 3 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B5").Value = 0.7964
$ws.Range("C5").Value = 0.4608
$ws.Range("D5").Value = 0.0839
$ws.Range("B6").Value = -4.1308
$ws.Range("C6").Value = 31132.3064
$ws.Range("D6").Value = 0.9999
$ws.Range("B7").Value = -0.1033
$ws.Range("C7").Value = 0.5976
$ws.Range("D7").Value = 0.8627
$ws.Range("B10").Value = 0.3199
$ws.Range("C10").Value = 0.4924
$ws.Range("D10").Value = 0.516

# --- Sheet 4 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B5").Value = 0.0062
$ws.Range("C5").Value = 0.0808
$ws.Range("D5").Value = 0.9385
$ws.Range("B6").Value = -0.0508
$ws.Range("C6").Value = 0.09
$ws.Range("D6").Value = 0.5723
$ws.Range("B7").Value = 0.0425
$ws.Range("C7").Value = 0.0822
$ws.Range("D7").Value = 0.6055
$ws.Range("B10").Value = -0.1305
$ws.Range("C10").Value = 0.0857
$ws.Range("D10").Value = 0.128

# --- Sheet 5 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B5").Value = 0.5066
$ws.Range("C5").Value = 0.4029
$ws.Range("D5").Value = 0.2086
$ws.Range("B6").Value = 0.2732
$ws.Range("C6").Value = 0.4071
$ws.Range("D6").Value = 0.5022
$ws.Range("B7").Value = 0.3417
$ws.Range("C7").Value = 0.4097
$ws.Range("D7").Value = 0.4043
$ws.Range("B10").Value = 0.0394
$ws.Range("C10").Value = 0.4375
$ws.Range("D10").Value = 0.9282

# --- Sheet 6 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("B5").Value = 0.2152
$ws.Range("C5").Value = 0.1585
$ws.Range("D5").Value = 0.1745
$ws.Range("B6").Value = -0.3065
$ws.Range("C6").Value = 0.2529
$ws.Range("D6").Value = 0.2254
$ws.Range("B7").Value = 0.2063
$ws.Range("C7").Value = 0.1625
$ws.Range("D7").Value = 0.2041
$ws.Range("B10").Value = 0.061
$ws.Range("C10").Value = 0.1622
$ws.Range("D10").Value = 0.707

# --- Sheet 7 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("B5").Value = 0.1365
$ws.Range("C5").Value = 0.6849
$ws.Range("D5").Value = 0.842
$ws.Range("B6").Value = 0.607
$ws.Range("C6").Value = 0.5845
$ws.Range("D6").Value = 0.299
$ws.Range("B7").Value = -0.2895
$ws.Range("C7").Value = 0.7234
$ws.Range("D7").Value = 0.6891
$ws.Range("B10").Value = -2.0377
$ws.Range("C10").Value = 2189.6613
$ws.Range("D10").Value = 0.9993

# --- Sheet 8 (): update coefficient table values ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("B5").Value = -0.0263
$ws.Range("C5").Value = 0.2367
$ws.Range("D5").Value = 0.9114
$ws.Range("B6").Value = -0.1329
$ws.Range("C6").Value = 0.2392
$ws.Range("D6").Value = 0.5783
$ws.Range("B7").Value = 0.0402
$ws.Range("C7").Value = 0.2417
$ws.Range("D7").Value = 0.8679
$ws.Range("B10").Value = -0.7947
$ws.Range("C10").Value = 0.2635
$ws.Range("D10").Value = 0.0026
